$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.649.22'
$ws.Range('E2').Value = '  -6.36%  '
$ws.Range('D3').Value = '3.766.50'
$ws.Range('E3').Value = '  -5.13%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.38%  '
$ws.Range('D7').Value = '3.761.64'
$ws.Range('E7').Value = '  -5.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.639'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.997'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('E10').Value = '  -6.37%  '
$ws.Range('E11').Value = '  -9.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.86%  '
$ws.Range('E13').Value = '  -9.74%  '
$ws.Range('E14').Value = '  -3.43%  '
$ws.Range('D15').Value = '4.353.48'
$ws.Range('E15').Value = '  -5.48%  '
$ws.Range('D16').Value = '3.758.34'
$ws.Range('E16').Value = '  -5.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.77%  '
$ws.Range('E19').Value = '  -6.32%  '
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').Value = '68.497.41'
$ws.Range('E21').Value = '  -6.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '414.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.41%  '
$ws.Range('E26').Value = '  -8.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  -7.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.91'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.32%  '
$ws.Range('E34').Value = '  -7.66%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '619.78'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '44.46'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.34%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0940'
$ws.Range('E37').Value = '  -10.69%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '65.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.06%  '
$ws.Range('E39').Value = '  -4.80%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.68%  '
$ws.Range('E43').Value = '  -5.56%  '
$ws.Range('E44').Value = '  -8.81%  '
$ws.Range('E45').Value = '  -7.18%  '
$ws.Range('E46').Value = '  +3.17%  '
$ws.Range('E47').Value = '  -9.48%  '
$ws.Range('E48').Value = '  -7.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.20%  '
$ws.Range('D50').Value = '2.780.71'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.27%  '
